# Applies the TPM-data update to the Pdgfc-Pdgfra LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01846133333333333
$ws.Range("H2").Value = 0.055384
$ws.Range("I2").Value = 0.005172740524168673
$ws.Range("J2").Value = 0.005172740524168674
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2662156666666667
$ws.Range("N2").Value = 0.7986470000000001
$ws.Range("O2").Value = 0.0009813702709097034
$ws.Range("P2").Value = 0.0009813702709097034
$ws.Range("Q2").Value = 0.00491469616088889
$ws.Range("R2").Value = 0.044232265448
$ws.Range("S2").Value = [double]"5.076373769549012E-06"
$ws.Range("T2").Value = [double]"5.076373769549013E-06"

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01846133333333333
$ws.Range("H3").Value = 0.055384
$ws.Range("I3").Value = 0.005172740524168673
$ws.Range("J3").Value = 0.005172740524168674
$ws.Range("O3").Value = 0.998256289001958
$ws.Range("P3").Value = 0.998256289001958
$ws.Range("Q3").Value = 4.999261233574223
$ws.Range("R3").Value = 44.99335110216801
$ws.Range("S3").Value = 0.005163720759626663
$ws.Range("T3").Value = 0.005163720759626664

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01846133333333333
$ws.Range("H4").Value = 0.055384
$ws.Range("I4").Value = 0.005172740524168673
$ws.Range("J4").Value = 0.005172740524168674
$ws.Range("M4").Value = 0.2067996666666667
$ws.Range("N4").Value = 0.620399
$ws.Range("O4").Value = 0.000762340727132399
$ws.Range("P4").Value = 0.0007623407271323989
$ws.Range("Q4").Value = 0.003817797579555556
$ws.Range("R4").Value = 0.034360178216
$ws.Range("S4").Value = [double]"3.943390772461973E-06"
$ws.Range("T4").Value = [double]"3.943390772461974E-06"

# Row 5
$ws.Range("I5").Value = 0.5495916323842472
$ws.Range("J5").Value = 0.5495916323842472
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2662156666666667
$ws.Range("N5").Value = 0.7986470000000001
$ws.Range("O5").Value = 0.0009813702709097034
$ws.Range("P5").Value = 0.0009813702709097034
$ws.Range("Q5").Value = 0.5221750198207777
$ws.Range("R5").Value = 4.699575178387001
$ws.Range("S5").Value = 0.0005393528891626348
$ws.Range("T5").Value = 0.0005393528891626348

# Row 6
$ws.Range("I6").Value = 0.5495916323842472
$ws.Range("J6").Value = 0.5495916323842472
$ws.Range("O6").Value = 0.998256289001958
$ws.Range("P6").Value = 0.998256289001958
$ws.Range("S6").Value = 0.5486333034104269
$ws.Range("T6").Value = 0.5486333034104269

# Row 7
$ws.Range("I7").Value = 0.5495916323842472
$ws.Range("J7").Value = 0.5495916323842472
$ws.Range("M7").Value = 0.2067996666666667
$ws.Range("N7").Value = 0.620399
$ws.Range("O7").Value = 0.000762340727132399
$ws.Range("P7").Value = 0.0007623407271323989
$ws.Range("Q7").Value = 0.4056321004421111
$ws.Range("R7").Value = 3.650688903979
$ws.Range("S7").Value = 0.0004189760846576892
$ws.Range("T7").Value = 0.0004189760846576891

# Row 8
$ws.Range("I8").Value = 0.445235627091584
$ws.Range("J8").Value = 0.445235627091584
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2662156666666667
$ws.Range("N8").Value = 0.7986470000000001
$ws.Range("O8").Value = 0.0009813702709097034
$ws.Range("P8").Value = 0.0009813702709097034
$ws.Range("Q8").Value = 0.4230248582804445
$ws.Range("R8").Value = 3.807223724524
$ws.Range("S8").Value = 0.0004369410079775195
$ws.Range("T8").Value = 0.0004369410079775195

# Row 9
$ws.Range("I9").Value = 0.445235627091584
$ws.Range("J9").Value = 0.445235627091584
$ws.Range("O9").Value = 0.998256289001958
$ws.Range("P9").Value = 0.998256289001958
$ws.Range("S9").Value = 0.4444592648319043
$ws.Range("T9").Value = 0.4444592648319043

# Row 10
$ws.Range("I10").Value = 0.445235627091584
$ws.Range("J10").Value = 0.445235627091584
$ws.Range("M10").Value = 0.2067996666666667
$ws.Range("N10").Value = 0.620399
$ws.Range("O10").Value = 0.000762340727132399
$ws.Range("P10").Value = 0.0007623407271323989
$ws.Range("Q10").Value = 0.3286110121897778
$ws.Range("R10").Value = 2.957499109708
$ws.Range("S10").Value = 0.0003394212517022478
$ws.Range("T10").Value = 0.0003394212517022478

Write-Output "Updated TPM values for rows 2-10"
